$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 12: PTLT account with password 12345 and UserType "super"
$ws.Range("A12").Value = "PTLT"
$ws.Range("B12").Value = 12345
$ws.Range("C12").Value = "super"

# Update the selected cell to D12 to match the saved selection state
$ws.Range("D12").Select()
